{"js": "const body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n\nconst replacements = [\n  [\"2025-06-19 Thursday\", \"2025-06-20 Friday\"],\n  [\"295\u00f77=\", \"680\u00f79=\"],\n  [\"221\u00f74=\", \"233\u00f73=\"],\n  [\"354\u00f74=\", \"305\u00f74=\"],\n  [\"951\u00f72=\", \"867\u00f74=\"],\n  [\"425\u00f72=\", \"473\u00f77=\"],\n  [\"935\u00f79=\", \"610\u00f72=\"],\n  [\"978\u00f79=\", \"763\u00f72=\"],\n  [\"976\u00f74=\", \"402\u00f77=\"],\n  [\"237\u00f77=\", \"522\u00f77=\"],\n  [\"686\u00f78=\", \"978\u00f78=\"],\n  [\"985\u00f78=\", \"270\u00f76=\"],\n  [\"200\u00f76=\", \"173\u00f78=\"],\n  [\"643\u00f72=\", \"389\u00f76=\"],\n  [\"841\u00f77=\", \"311\u00f72=\"],\n  [\"673\u00f72=\", \"624\u00f77=\"],\n  [\"524\u00f76=\", \"623\u00f75=\"],\n  [\"876\u00f75=\", \"427\u00f77=\"],\n  [\"850\u00f74=\", \"513\u00f74=\"],\n  [\"838\u00f79=\", \"346\u00f77=\"],\n  [\"492\u00f74=\", \"577\u00f75=\"],\n  [\"923\u00f76=\", \"127\u00f72=\"],\n  [\"912\u00f78=\", \"169\u00f74=\"],\n  [\"842\u00f72=\", \"781\u00f77=\"],\n  [\"940\u00f76=\", \"762\u00f73=\"],\n  [\"305\u00f73=\", \"673\u00f76=\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-19 Thursday\", \"2025-06-20 Friday\"),\n    @(\"295\u00f77=\", \"680\u00f79=\"),\n    @(\"221\u00f74=\", \"233\u00f73=\"),\n    @(\"354\u00f74=\", \"305\u00f74=\"),\n    @(\"951\u00f72=\", \"867\u00f74=\"),\n    @(\"425\u00f72=\", \"473\u00f77=\"),\n    @(\"935\u00f79=\", \"610\u00f72=\"),\n    @(\"978\u00f79=\", \"763\u00f72=\"),\n    @(\"976\u00f74=\", \"402\u00f77=\"),\n    @(\"237\u00f77=\", \"522\u00f77=\"),\n    @(\"686\u00f78=\", \"978\u00f78=\"),\n    @(\"985\u00f78=\", \"270\u00f76=\"),\n    @(\"200\u00f76=\", \"173\u00f78=\"),\n    @(\"643\u00f72=\", \"389\u00f76=\"),\n    @(\"841\u00f77=\", \"311\u00f72=\"),\n    @(\"673\u00f72=\", \"624\u00f77=\"),\n    @(\"524\u00f76=\", \"623\u00f75=\"),\n    @(\"876\u00f75=\", \"427\u00f77=\"),\n    @(\"850\u00f74=\", \"513\u00f74=\"),\n    @(\"838\u00f79=\", \"346\u00f77=\"),\n    @(\"492\u00f74=\", \"577\u00f75=\"),\n    @(\"923\u00f76=\", \"127\u00f72=\"),\n    @(\"912\u00f78=\", \"169\u00f74=\"),\n    @(\"842\u00f72=\", \"781\u00f77=\"),\n    @(\"940\u00f76=\", \"762\u00f73=\"),\n    @(\"305\u00f73=\", \"673\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
